# Updated cryptos list on Fri Apr 19 05:54:30 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Price/Coin/Link cells in this sheet are plain text (no number format).
    # Force text storage so numeric-looking strings (e.g. "544.79") don't
    # get reinterpreted as numbers by Excel's usual type inference, then
    # restore the default (unstyled) cell style so no stray style index is
    # introduced.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# --- Simple per-row Price (D) / Volume(1h) (E) updates ---

Set-TextCell "D2" "62.272.11"
$ws.Range("E2").Value = "  +1.41%  "

Set-TextCell "D3" "3.011.15"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  -0.08%  "

Set-TextCell "D5" "544.79"
$ws.Range("E5").Value = "  -0.67%  "

Set-TextCell "D6" "139.42"
$ws.Range("E6").Value = "  +3.69%  "

Set-TextCell "D7" "0.999"
$ws.Range("E7").Value = "  -0.09%  "

Set-TextCell "D8" "3.006.04"
$ws.Range("E8").Value = "  +0.29%  "

Set-TextCell "D9" "0.489"
$ws.Range("E9").Value = "  -1.49%  "

Set-TextCell "D10" "6.69"
$ws.Range("E10").Value = "  +10.52%  "

$ws.Range("E11").Value = "  -0.04%  "

Set-TextCell "D12" "0.445"
$ws.Range("E12").Value = "  -1.35%  "

Set-TextCell "D13" "0.0000221"
$ws.Range("E13").Value = "  -0.24%  "

Set-TextCell "D14" "33.97"
$ws.Range("E14").Value = "  -1.69%  "

Set-TextCell "D15" "3.487.11"
$ws.Range("E15").Value = "  -0.25%  "

Set-TextCell "D16" "62.292.97"
$ws.Range("E16").Value = "  +1.21%  "

Set-TextCell "D17" "3.007.56"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("E18").Value = "  -2.49%  "

Set-TextCell "D19" "6.56"
$ws.Range("E19").Value = "  -1.76%  "

Set-TextCell "D20" "466.61"
$ws.Range("E20").Value = "  -1.47%  "

Set-TextCell "D21" "13.36"
$ws.Range("E21").Value = "  +0.57%  "

Set-TextCell "D22" "0.654"
$ws.Range("E22").Value = "  -3.20%  "

Set-TextCell "D23" "7.22"
$ws.Range("E23").Value = "  +2.50%  "

# --- Rows 24 & 25 swap: Litecoin <-> InternetComputer(DFINITY), with new data ---

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D24" "12.62"
$ws.Range("E24").Value = "  +4.19%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "79.25"
$ws.Range("E25").Value = "  -1.02%  "

Set-TextCell "D26" "1.00"
$ws.Range("E26").Value = "  +0.10%  "

Set-TextCell "D27" "2.71"
$ws.Range("E27").Value = "  -0.07%  "

Set-TextCell "D28" "7.62"
$ws.Range("E28").Value = "  -2.60%  "

$ws.Range("E29").Value = "  +4.91%  "

$ws.Range("E30").Value = "  -0.17%  "

Set-TextCell "D31" "25.46"
$ws.Range("E31").Value = "  -0.94%  "

$ws.Range("E32").Value = "  -1.51%  "

Set-TextCell "D33" "2.33"
$ws.Range("E33").Value = "  +1.31%  "

Set-TextCell "D34" "5.55"
$ws.Range("E34").Value = "  +0.04%  "

Set-TextCell "D35" "54.68"
$ws.Range("E35").Value = "  -1.25%  "

Set-TextCell "D36" "5.84"
$ws.Range("E36").Value = "  -1.20%  "

Set-TextCell "D37" "450.97"
$ws.Range("E37").Value = "  -0.74%  "

Set-TextCell "D38" "0.0807"
$ws.Range("E38").Value = "  +0.94%  "

Set-TextCell "D39" "0.0390"
$ws.Range("E39").Value = "  +1.89%  "

Set-TextCell "D40" "2.940.54"
$ws.Range("E40").Value = "  -7.81%  "

Set-TextCell "D41" "0.115"
$ws.Range("E41").Value = "  -2.59%  "

Set-TextCell "D42" "8.09"
$ws.Range("E42").Value = "  -0.93%  "

Set-TextCell "D43" "2.59"
$ws.Range("E43").Value = "  +6.03%  "

Set-TextCell "D44" "26.80"
$ws.Range("E44").Value = "  +2.17%  "

Set-TextCell "D46" "0.247"
$ws.Range("E46").Value = "  +0.86%  "

# --- Rows 47 & 48 swap: Stellar <-> Fetch.AI, with new data ---

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D47" "2.00"
$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D48" "0.109"
$ws.Range("E48").Value = "  +0.69%  "

Set-TextCell "D49" "114.82"
$ws.Range("E49").Value = "  -2.77%  "

$subscript3 = [char]8323
Set-TextCell "D50" "0.0${subscript3}0498"
$ws.Range("E50").Value = "  +1.21%  "

$ws.Range("E51").Value = "  -3.20%  "
